# EnterpriseSupportDatasheet.pptx - footnote punctuation fixes on the
# "Support Coverage Hours" table (slide 4, graphicFrame id=25 "Table 6").
#
#   1. "Language support is only available in English and Japanese " (+trailing
#      space) -> "Language support is only available in English and Japanese."
#   2. "P2, P3, P4 cases are limited to business hours only in Japan" ->
#      "P2, P3, P4 cases are limited to business hours only in Japan."

$p = $ppt.ActivePresentation

# The slide with sldId=261 is the 4th slide in the deck.
$s = $p.Slides.Item(4)

# Locate the table shape by its persisted shape id (25, "Table 6") rather
# than a hard-coded collection index.
$tableShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Id -eq 25) {
        $tableShape = $candidate
        break
    }
}

$tbl = $tableShape.Table

# Row 3 is the merged (gridSpan=4) footnote row spanning all 4 columns; any
# column index reaches the same underlying cell/shape.
$cell = $tbl.Cell(3, 1)
$tr = $cell.Shape.TextFrame.TextRange

# --- Edit 1: tidy up the end of the "Language support ..." sentence ---
$full = $tr.Text
$marker = "Japanese"
$markerStart = $full.IndexOf($marker)
$afterMarkerPos = $markerStart + $marker.Length + 1   # 1-based char right after "Japanese"
$trailingChar = $tr.Characters($afterMarkerPos, 1)
$trailingChar.Text = "."

# --- Edit 2: add the missing period after "... only in Japan" ---
$full = $tr.Text
$lastPos = $full.Length
$lastChar = $tr.Characters($lastPos, 1)
$lastChar.Text = $lastChar.Text + "."
